# 1b - Design Decisions: "Final draft of this" edit
#
# Two changes:
#   1) Merge "filtered through" and " after being broken down through
#      normalization and convolution." into a single run (no visible text
#      change - just tidies up the run split left over from an earlier
#      edit).
#   2) Insert the word "were " before "defined using summations" (turning
#      "...computational building blocks defined using summations..." into
#      "...computational building blocks were defined using summations...").
#      Word's automatic "_GoBack" bookmark (which always marks the site of
#      the most recent edit) needs to be relocated to sit at that new typing
#      position, right before "defined".

$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# Re-save the sentence as a single run so "filtered through" and the
# trailing clause are no longer split across two <w:r> elements.
$range1 = $d.Content
$found1 = $range1.Find.Execute(
    "filtered through after being broken down through normalization and convolution.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "filtered through after being broken down through normalization and convolution.",
    2)

# --- Change 2 -------------------------------------------------------------
# Locate "defined using summations" and type "were " immediately before it.
$range2 = $d.Content
$found2 = $range2.Find.Execute(
    "defined using summations", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

$insertPoint = $d.Range($range2.Start, $range2.Start)
$insertPoint.InsertBefore("were ")

# Move the "_GoBack" bookmark to the new edit location (right before
# "defined", after the newly-typed "were ").
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
}

$range3 = $d.Content
$found3 = $range3.Find.Execute(
    "defined using summations", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
$bookmarkRange = $d.Range($range3.Start, $range3.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
